# Deploy updated output folder
# - Metadata!B8 "Date" value bumped to the new generation timestamp.
# - Elements sheet: the "Address Extension" binding strength/value set
#   entries (row 6) changed from required/ng-wards to example/nigeria-wards.
# - Column Z ("Binding Value Set") widens to fit the longer URL (bestFit).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-04T07:50:29+01:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("X6").Value = "example"
$elements.Range("Z6").Value = "https://nphcda.gov.ng/immunizationIG/ValueSet/nigeria-wards"

# Re-run the sheet's "best fit" column-width pass on column Z (26) now that
# its text is longer.
$elements.Columns.Item(26).ColumnWidth = 49.5

Write-Output "Updated Date, binding strength/value set, and column Z width"
